$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D (Price) and E (Volume(1h)) hold numeric-looking data stored as TEXT
# (inline strings) in the source sheet. Forcing NumberFormat to "@" before
# assignment keeps Excel from re-typing the write as a number/percentage,
# then resetting the style back to "Normal" restores the original (default)
# cell styling so only the text content changes.
$changes = @(
    @{ Row = 2; Col = 4; Value = "304.38" },
    @{ Row = 2; Col = 5; Value = "0.77%" },
    @{ Row = 3; Col = 4; Value = "35.65" },
    @{ Row = 3; Col = 5; Value = "-4.72%" },
    @{ Row = 4; Col = 4; Value = "5.083" },
    @{ Row = 4; Col = 5; Value = "1.62%" },
    @{ Row = 5; Col = 5; Value = "0.50%" },
    @{ Row = 6; Col = 4; Value = "2.127" },
    @{ Row = 6; Col = 5; Value = "-2.91%" },
    @{ Row = 7; Col = 4; Value = "7.908" },
    @{ Row = 8; Col = 4; Value = "0.9193" },
    @{ Row = 8; Col = 5; Value = "0.57%" },
    @{ Row = 9; Col = 4; Value = "0.09747" },
    @{ Row = 9; Col = 5; Value = "0.04%" },
    @{ Row = 10; Col = 4; Value = "0.1859" },
    @{ Row = 10; Col = 5; Value = "-1.77%" },
    @{ Row = 11; Col = 4; Value = "0.08581" },
    @{ Row = 11; Col = 5; Value = "-0.09%" },
    @{ Row = 12; Col = 4; Value = "0.03553" },
    @{ Row = 12; Col = 5; Value = "0.75%" },
    @{ Row = 13; Col = 4; Value = "0.09949" },
    @{ Row = 13; Col = 5; Value = "-0.15%" },
    @{ Row = 14; Col = 4; Value = "0.001437" },
    @{ Row = 14; Col = 5; Value = "-3.55%" },
    @{ Row = 15; Col = 4; Value = "0.005632" },
    @{ Row = 15; Col = 5; Value = "-1.71%" },
    @{ Row = 16; Col = 4; Value = "3.467" },
    @{ Row = 16; Col = 5; Value = "0.09%" },
    @{ Row = 17; Col = 4; Value = "4.095" },
    @{ Row = 17; Col = 5; Value = "1.40%" },
    @{ Row = 18; Col = 4; Value = "2.556" },
    @{ Row = 18; Col = 5; Value = "22.75%" },
    @{ Row = 19; Col = 4; Value = "0.3425" },
    @{ Row = 19; Col = 5; Value = "-1.10%" },
    @{ Row = 20; Col = 4; Value = "5.226" },
    @{ Row = 20; Col = 5; Value = "9.68%" },
    @{ Row = 21; Col = 4; Value = "0.1309" },
    @{ Row = 21; Col = 5; Value = "0.55%" },
    @{ Row = 22; Col = 5; Value = "-0.07%" },
    @{ Row = 23; Col = 4; Value = "0.04546" },
    @{ Row = 23; Col = 5; Value = "-2.07%" },
    @{ Row = 24; Col = 4; Value = "0.005060" },
    @{ Row = 24; Col = 5; Value = "5.52%" },
    @{ Row = 25; Col = 4; Value = "0.001235" },
    @{ Row = 25; Col = 5; Value = "0.41%" },
    @{ Row = 27; Col = 4; Value = "0.0004752" },
    @{ Row = 27; Col = 5; Value = "0.02%" },
    @{ Row = 39; Col = 4; Value = "0.01843" },
    @{ Row = 39; Col = 5; Value = "4.82%" },
    @{ Row = 40; Col = 4; Value = "0.04723" },
    @{ Row = 40; Col = 5; Value = "-0.07%" },
    @{ Row = 41; Col = 4; Value = "0.007491" },
    @{ Row = 41; Col = 5; Value = "-7.11%" },
    @{ Row = 42; Col = 4; Value = "0.1399" },
    @{ Row = 42; Col = 5; Value = "0.43%" },
    @{ Row = 43; Col = 4; Value = "0.007747" },
    @{ Row = 44; Col = 4; Value = "0.002205" },
    @{ Row = 44; Col = 5; Value = "2.00%" },
    @{ Row = 45; Col = 5; Value = "5.70%" },
    @{ Row = 46; Col = 4; Value = "0.00006322" },
    @{ Row = 46; Col = 5; Value = "5.87%" },
    @{ Row = 47; Col = 4; Value = "0.00000000750" },
    @{ Row = 47; Col = 5; Value = "-0.02%" },
    @{ Row = 48; Col = 4; Value = "0.0005799" },
    @{ Row = 48; Col = 5; Value = "-0.04%" },
    @{ Row = 49; Col = 4; Value = "46.60" },
    @{ Row = 49; Col = 5; Value = "509.93%" },
    @{ Row = 50; Col = 5; Value = "-25.63%" },
    @{ Row = 51; Col = 4; Value = "0.00002101" },
    @{ Row = 51; Col = 5; Value = "-0.02%" }
)

foreach ($chg in $changes) {
    $cell = $ws.Cells.Item($chg.Row, $chg.Col)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
    $cell.Style = "Normal"
}
